$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared string: "滚动条样式" will be introduced via A12's text value ---

# --- Row 7: add C7 ("check mark") and D7 (date 2017-02-25), matching the
#     style already used by the same columns in row 8 ---
$ws.Range("C8").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = $ws.Range("C8").Value()

$ws.Range("D8").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 42791

# --- Row 12 (new row): "滚动条样式" task, proposed 2017-02-25, done, resolved 2017-02-25 ---
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "滚动条样式"

$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 42791

$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = $ws.Range("C11").Value()

$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 42791

$ws.Rows.Item(12).RowHeight = 24.75

# --- Update the active selection to F8, matching the post-edit cursor position ---
[void]$ws.Range("F8").Select()
